# BarGuestDB.xlsx update: hook up the real ScriptableObject guest data
# (Guests / Scripts / CocktailProblems sheets) replacing the old
# placeholder "베일 존 / 임주원" test rows with "미아 / 다니엘" content,
# and dropping the extra "지혁" row from CocktailProblems.

$wb = $excel.ActiveWorkbook

$wsGuests = $wb.Worksheets.Item("Guests")
$wsScripts = $wb.Worksheets.Item("Scripts")
$wsCocktail = $wb.Worksheets.Item("CocktailProblems")

# ---------------------------------------------------------------------
# Scripts sheet: replace the character dialogue rows.
# guest_code 0001 used to be "베일 존" (6 lines) - now "미아" (8 lines).
# guest_code 0002 used to be "임주원" (6 lines) - now "다니엘" (12 lines).
# ---------------------------------------------------------------------
$scriptsRows = @(
    @("0001", 0, "미아", "여긴 본 적이 없는데 최근에 개업하셨나봐요?"),
    @("0001", 0, "미아", "이런 시골 촌에도 입주를 하는 사람이 있구나~"),
    @("0001", 0, "미아", "뭐 칵테일을 먹어본 적은 없지만.."),
    @("0001", 0, "미아", "나는 상큼하고 신 걸 좋아해요."),
    @("0001", 0, "미아", "뭐 주문 내용은 이해하셨겠죠?"),
    @("0001", 1, "미아", "내가 칵테일에 대해 아는 건 없지만.."),
    @("0001", 1, "미아", "옆 집 펍보다는 그나마 나은 수준이네요."),
    @("0001", 1, "미아", "뭐 현실적으로 번창하라는 말은 못하겠지만.. 열심히 해보세요~"),
    @("0002", 0, "다니엘", "여기가 친구가 말한 바가 맞겠지?"),
    @("0002", 0, "다니엘", "새로 생긴 바가 있다고 해서 친구놈이 실험쥐로 날 먼저 보냈어"),
    @("0002", 0, "다니엘", "이 동네 사람들은 바라는 곳을 한 번도 가본 적이 없어서 아마 시내에 한 번이라도 가본 나를 앞세운 거겠지"),
    @("0002", 0, "다니엘", "뭐 나도 사실상 바를 한 번도 가본 적은 없지만 말이야 하.."),
    @("0002", 0, "다니엘", "하하.."),
    @("0002", 0, "다니엘", "너 생각보다 말이 없구나?"),
    @("0002", 0, "다니엘", "됐어 그럼. 내가 지금 외롭기도 하고.. 우울감도 좀 있는 것 같아서 아무한테나 말을 걸고 싶었거든"),
    @("0002", 0, "다니엘", "아무튼 바텐더, 오늘 내 기분대로 한 잔~"),
    @("0002", 0, "다니엘", "이렇게 주문하는게 맞나?"),
    @("0002", 1, "다니엘", "뭐 맥주만 마시던 내가 원하던 느낌은 아니지만.."),
    @("0002", 1, "다니엘", "일단 취하니까 기분은 좋네~"),
    @("0002", 1, "다니엘", "번창하쇼~ 아니 번창은 힘들겠구나..")
)

$r = 2
foreach ($row in $scriptsRows) {
    $wsScripts.Cells.Item($r, 1).Value = $row[0]
    $wsScripts.Cells.Item($r, 2).Value = $row[1]
    $wsScripts.Cells.Item($r, 3).Value = $row[2]
    $wsScripts.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Widen the script column to fit the new, longer lines.
$wsScripts.Columns.Item(4).ColumnWidth = 89.5

# ---------------------------------------------------------------------
# CocktailProblems sheet: drop the extra "아니요" row that belonged to
# the old 4-guest placeholder roster (guest_code 0001, scan_type 1).
# ---------------------------------------------------------------------
$wsCocktail.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# View state: Guests becomes the active sheet/tab with D6 selected;
# Scripts keeps A19 selected; CocktailProblems keeps E4 selected.
# ---------------------------------------------------------------------
$wsScripts.Activate()
$wsScripts.Range("A19").Select()

$wsCocktail.Activate()
$wsCocktail.Range("E4").Select()

$wsGuests.Activate()
$wsGuests.Range("D6").Select()

Write-Output "BarGuestDB updated"
